$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A43").Value = 10041
$ws.Range("B43").Value = "message_10041_physical_address_record_updated_successfully"
$ws.Range("D43").Value = "Success"

$ws.Range("A44").Value = 10042
$ws.Range("B44").Value = "message_10042_physical_address_record_deleted_successfully"
$ws.Range("D44").Value = "Success"

$ws.Range("A45").Select() | Out-Null
